$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 0
$wsExhibition.Range("F3").Value = 0
$wsExhibition.Range("F4").Value = 1639
$wsExhibition.Range("F5").Value = 0
$wsExhibition.Range("F6").Value = 0
$wsExhibition.Range("F8").Value = 0

# Sheet "全部类型" (All Types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 394
$wsAll.Range("F5").Value = 0
$wsAll.Range("F7").Value = 418
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 0
